# Duplicate the "LKT 8HED3" sheet to a new sheet "LKT 8HED3A" placed at the
# end of the workbook, and make the new sheet the active tab/sheet.

$wb = $excel.ActiveWorkbook

$source = $wb.Worksheets.Item("LKT 8HED3")

# Copy the source sheet to after the last sheet in the workbook.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$source.Copy($null, $lastSheet)

# The newly created sheet is now the last sheet and is active; rename it.
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "LKT 8HED3A"

# Make sure the new sheet is the active/selected sheet.
$newSheet.Activate()
$newSheet.Select()
